# clickEvent() / save(): persist the current "info" record and append the
# line items that make up the saving tree onto their respective sheets.

$wb = $excel.ActiveWorkbook

# --- "info" sheet: the saved record's summary fields (kept as text) ---
$infoSheet = $wb.Worksheets.Item("info")
$infoValues = @("123", "123", "123", "2")
for ($col = 1; $col -le $infoValues.Length; $col++) {
    $cell = $infoSheet.Cells.Item(1, $col)
    # Force text storage so numeric-looking strings stay strings, then drop
    # the temporary number format so no extra style sticks to the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $infoValues[$col - 1]
    $cell.ClearFormats()
}

# --- "items" sheet: saved shopping/expense tree (name, unit, price, qty, total) ---
$itemsSheet = $wb.Worksheets.Item("items")

$items = @(
    @("접시140",   "개", 1000, 1, 1000),
    @("신라면(컵)", "개", 930,  1, 930),
    @("신라면(컵)", "개", 930,  1, 930),
    @("식탁보",     "개", 6000, 1, 6000)
)

for ($row = 1; $row -le $items.Length; $row++) {
    $item = $items[$row - 1]

    $nameCell = $itemsSheet.Cells.Item($row, 1)
    $nameCell.NumberFormat = "@"
    $nameCell.Value = $item[0]
    $nameCell.ClearFormats()

    $unitCell = $itemsSheet.Cells.Item($row, 2)
    $unitCell.NumberFormat = "@"
    $unitCell.Value = $item[1]
    $unitCell.ClearFormats()

    $itemsSheet.Cells.Item($row, 3).Value = $item[2]
    $itemsSheet.Cells.Item($row, 4).Value = $item[3]
    $itemsSheet.Cells.Item($row, 5).Value = $item[4]
}
